$wb = $excel.ActiveWorkbook

# Both the "展览" (exhibition) sheet and the "全部类型" (all types) sheet
# contain the same data table and both need the same updated values.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1700
    $ws.Range("F4").Value = 331
    $ws.Range("F6").Value = 779
    $ws.Range("F8").Value = 5812
}
